$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the new test name.
$ws.Name = "AddCustomerTest"

# Header row (row 1) + data row (row 2), filled in the same order the
# original test authored them (lastName/firstName/postCode, then the
# sample row, then the alertText column added afterwards).
$ws.Range("B1").Value = "lastName"
$ws.Range("A1").Value = "firstName"
$ws.Range("C1").Value = "postCode"

$ws.Range("A2").Value = "Sohaib"
$ws.Range("B2").Value = "Majeed"
$ws.Range("C2").Value = "123wp"

$ws.Range("D1").Value = "alertText"
$ws.Range("D2").Value = "Customer added successfully"

# Header row is bold.
$ws.Range("A1:D1").Font.Bold = $true

# Column widths sized to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 9
$ws.Columns.Item(2).ColumnWidth = 8.666666666666668
$ws.Columns.Item(3).ColumnWidth = 8.666666666666668
$ws.Columns.Item(4).ColumnWidth = 42.16666666666667

# Print orientation.
$ws.PageSetup.Orientation = 1

# Leave the cursor where the user left it after entering the data.
[void]$ws.Range("D3").Select()
